$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 54 is no longer the last data row: drop its bottom border
# (matches the border-8 family already used on rows 47-53).
$ws.Range("A54:E54").Borders(9).LineStyle = -4142

# Seed row 55 formatting from cells that already carry the exact
# target styles (full border box, the "last row" look).
$ws.Range("C2").Copy()
$ws.Range("A55").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C55").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B55").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D55").PasteSpecial(-4122)
$ws.Range("E47").Copy()
$ws.Range("E55").PasteSpecial(-4122)

# New knowledge-base entry: ACN customer support
$ws.Range("A55").Value = 'ACN 고객 지원하기'
$ws.Range("B55").Value = 45457
$ws.Range("C55").Value = 'CP00013'
$ws.Range("D55").Value = 'ACN 고객 지원 방법에 대해 알아봅니다.
<br>ACN은 <strong>Automatic Crash Notification</strong>의 약자로 한국어로 자동 충격(충돌) 알림을 의미합니다.
<br>
<br>ACN 통화의 경우 CCC에서 <a style="color: red; font-weight: bold;">ACN</a>으로 표기되며, Saleforce 내 Genesys에서도 <a style="color: red; font-weight: bold;">ACN</a>으로 표기됩니다.
<br>
<br.아래에 따라 진행합니다.
<br><em style="font-weight: bold; font-size: 17pt;">볼보 긴급 지원 서비스 입니다. 고객님 괜찮으십니까?"</em>
<br>
<br>응답이 있는 경우
<br><em style="font-weight: bold; font-size: 17pt;">"고객님 구조 요청이 필요하십니까?"</em>
<br>필요한 경우
<br><em style="font-weight: bold; font-size: 17pt;">"현재 OOOO위치로 확인됩니다. 해당 위치로 구조 요청 진행해도 될까요? 고객님 성함과 연락처 말씀 부탁드립니다."</em>
<br><strong>E-Call 유지 상태에서 119에 신고한 후 3자 통화로 전환하여 확인합니다.</strong>
<br><a style="font-weight: bold;">*119에 신고하는 경우 고객 위치 (시)의 소방으로 연결해달라고 설명합니다.</a>
<br><a style="font-wieght: bold; color: red;">119에 신고하는 경우: </a>
<br><em style="font-weight: bold; font-size: 17pt;">1. "안녕하십니까 볼보고객지원팀입니다. OO시로 연결 부탁드립니다."
<br>2. (연결 된 경우) "안녕하십니까 볼보고객지원팀입니다. 현재 고객 사고로 인해 신고 접수 드립니다. 현재 위치 알려드리겠습니다. (CCC 위치 확인 후 전달)"</em>
<br>3. 신고 접수 후 119와 통화가 끝어지면 아래 단계를 수행합니다.
<br><em style="font-weight: bold; font-size: 17pt;">3-1 "구조대가 올 때까지 통화를 유지해도 되겠습니까? 그렇지 않은 경우 연락처로 15분 후에 연락드리겠습니다."</em>
<br>3-2 고객이 콜백을 요청한 경우
<br>- 15분 후 고객에게 연락합니다.
<br>
<br>고객이 통화 유지를 요청한 경우
<br>- 통화를 유지합니다. / 사이렌 소리가 들리면 고객에게 통화 종료 의사를 묻고 종료합니다.
<br>
<br><em style="font-weight: bold; font-size: 17pt;">119 신고를 원치 않는 경우:</em>
<br>긴급 출동 서비스 지원에 대해 안내 합니다.
<br>
<br>
응답이 없는 경우
<em style="font-weight: bold; font-size: 17pt;">"고객님 들리십니까?"</em> 지속 반복 / 30초 이상 무응답 시 아래 단계 수행
즉시 119에 신고합니다. / 3자 통화로 전환합니다.
1. 119에 신고하는 경우: <em style="font-weight: bold; font-size: 17pt;">"안녕하십니까 볼보고객지원팀입니다. OO시로 연결 부탁드립니다."</em>
2. <em style="font-weight: bold; font-size: 17pt;">"저희 고객께서 현재 차량 충돌로 자동 신고가 접수되어 신고 드립니다. 현재 위치 알려드리겠습니다. (CCC 위치 확인 후 전달)"</em>
3. 119와 통화가 끝어지면 고객을 계속 부릅니다. 통화 상에서 사이렌이 울리고 구조대가 도착할 때까지 유지합니다. 구조대가 도착하면 통화 중인 상태를 알립니다.
4. 고객이 응답이 없는 경우 구조대를 통해 통화를 끊어도 될지 물어봅니다.'
$ws.Range("E55").Value = 'RSA'

# Row heights: row 54 keeps its height, row 55 is a very long entry
# that hits Excel's maximum row height.
$ws.Rows.Item(55).RowHeight = 409.6

# Selection / view mirrors the author's final scroll position.
$win = $excel.ActiveWindow
$win.ScrollRow = 54
$win.ScrollColumn = 4
$ws.Range("H55").Select()

